$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.580.75"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").Value = "2.630.14"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "604.87"
$ws.Range("E5").Value = "  +1.51%  "

$ws.Range("D6").Value = "154.73"
$ws.Range("E6").Value = "  -0.27%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "0.547"
$ws.Range("E8").Value = "  +0.81%  "

$ws.Range("D9").Value = "2.628.41"
$ws.Range("E9").Value = "  +0.22%  "

$ws.Range("E10").Value = "  +6.70%  "

$ws.Range("E11").Value = "  +0.53%  "

$ws.Range("D12").Value = "5.23"
$ws.Range("E12").Value = "  -0.18%  "

$ws.Range("D13").Value = "0.353"
$ws.Range("E13").Value = "  -1.55%  "

$ws.Range("E14").Value = "  +0.03%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.115.57"
$ws.Range("E15").Value = "  +0.71%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000184"
$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("D17").Value = "67.714.33"
$ws.Range("E17").Value = "  +0.61%  "

$ws.Range("D18").Value = "2.636.09"
$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "11.30"
$ws.Range("E19").Value = "  -0.85%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "365.94"
$ws.Range("E20").Value = "  +0.22%  "

$ws.Range("D21").Value = "7.60"
$ws.Range("E21").Value = "  -3.37%  "

$ws.Range("E22").Value = "  -0.54%  "

$ws.Range("D23").Value = "2.13"
$ws.Range("E23").Value = "  +4.99%  "

$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").Value = "70.53"
$ws.Range("E25").Value = "  +3.15%  "

$ws.Range("D26").Value = "10.06"
$ws.Range("E26").Value = "  -2.16%  "

$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.764.06"
$ws.Range("E27").Value = "  +0.14%  "

$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0000104"
$ws.Range("E28").Value = "  -1.40%  "

$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").Value = "585.89"
$ws.Range("E29").Value = "  -5.12%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.79%  "

$ws.Range("D31").Value = "1.43"
$ws.Range("E31").Value = "  -2.61%  "

$ws.Range("D32").Value = "7.91"
$ws.Range("E32").Value = "  -2.48%  "

$ws.Range("E33").Value = "  -0.63%  "

$ws.Range("D34").Value = "0.130"
$ws.Range("E34").Value = "  -3.04%  "

$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("D36").Value = "1.54"
$ws.Range("E36").Value = "  -3.09%  "

$ws.Range("D37").Value = "4.97"
$ws.Range("E37").Value = "  -1.07%  "

$ws.Range("D38").Value = "19.52"
$ws.Range("E38").Value = "  -0.27%  "

$ws.Range("D39").Value = "157.60"
$ws.Range("E39").Value = "  +2.22%  "

$ws.Range("E40").Value = "  -0.23%  "

$ws.Range("D41").Value = "5.33"
$ws.Range("E41").Value = "  -3.61%  "

$ws.Range("D42").Value = "1.84"
$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("D43").Value = "2.62"
$ws.Range("E43").Value = "  -2.29%  "

$ws.Range("D44").Value = "41.22"
$ws.Range("E44").Value = "  -0.56%  "

$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "16.37"
$ws.Range("E46").Value = "  -0.67%  "

$ws.Range("D47").Value = "156.32"
$ws.Range("E47").Value = "  -0.10%  "

$ws.Range("D48").Value = "0.0₆0289"
$ws.Range("E48").Value = "  -6.00%  "

$ws.Range("D49").Value = "3.74"
$ws.Range("E49").Value = "  -1.23%  "

$ws.Range("D50").Value = "20.97"
$ws.Range("E50").Value = "  -1.85%  "

$ws.Range("D51").Value = "0.625"
$ws.Range("E51").Value = "  -0.42%  "
